# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# Update the OFF sheet (Target Depth Data - offense)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 498
$wsOff.Range("C2").Value = 348
$wsOff.Range("D2").Value = 107
$wsOff.Range("E2").Value = 55

# Update the DEF sheet (Target Depth Data - defense)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 421
$wsDef.Range("C2").Value = 290
$wsDef.Range("D2").Value = 95
$wsDef.Range("E2").Value = 38
$wsDef.Range("F2").Value = 6
